# Reorder the comma-separated names within the "Recorded By" (column G) cells.
# This mirrors the upstream sync: entries are re-joined in a different order
# without changing the underlying set of names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value2 = "system, backup@backdoor.com, System"
    }
}
